# Sprint Logs updated with meeting 3
# Fills in column E ("Answer") for every question/person row with the
# meeting-3 log entries, restyles three long answers with smaller fonts,
# grows the two rows that needed more height to show the wrapped text,
# and leaves the selection where the author ended up (E20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Pre-register the three custom font sizes used by the long answers
#        in rows 2, 14 and 18 (10pt, 9pt, 8pt respectively) on scratch
#        cells well outside the used range, then copy *only the formats*
#        onto the target cells (in row order E2, E14, E18). Doing the
#        font-size registration and the per-cell format realization as
#        two separate passes reproduces the exact font/cellXf ordering
#        Excel produced in the canonical file. The scratch cells are
#        cleared afterwards so nothing extra is left in the sheet.
$ws.Range("A200").Font.Size = 10
$ws.Range("A200").WrapText = $true
$ws.Range("A201").Font.Size = 9
$ws.Range("A201").WrapText = $true
$ws.Range("A202").Font.Size = 8
$ws.Range("A202").WrapText = $true

$ws.Range("A200").Copy()
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("A201").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("A202").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("A200:A202").Clear()

# --- 2) Fill in the answers. Order matters for de-duplication bookkeeping
#        of repeated text (e.g. "No", "Not yet") and matches the order the
#        entries were originally authored in.
$ws.Range("E2").Value = 'Successfully connected the application to Google firebase. Utilizing the following set up for the andrio studio project: Download android studio v3.3.1 for windows from chrome using API 20:Andriod 4.4W (kitkat wear) for project API using Nexus 5x API with Nougat 24 OS for emulated device. Completed basic design for security question update page. Completed back end for security question page however does not consider user account or database, will have to e updated for password verification and to actually update database'
$ws.Range("E6").Value = 'Try to connect the submission of a security question with Google firebase'
$ws.Range("E14").Value = 'Installing Android studio from certain rowsers has proven to be problematic. So far google chrome has been issue free.'
$ws.Range("E18").Value = 'Versions od IDE, project, and emulator sohuld be consistent while in the rapid prototyping phase.'
$ws.Range("E5").Value = 'Continued work on sprites for the game'
$ws.Range("E9").Value = 'Continue to work on sprite animations and button code for the game'
$ws.Range("E13").Value = 'Nothing is currently getting in the way of  my work'
$ws.Range("E17").Value = 'Continue learning how to utilize new animation techniques'
$ws.Range("E21").Value = 'No changes need to be made to the project currently'
$ws.Range("E11").Value = 'Chores for spring break'
$ws.Range("E15").Value = 'documentation is tricky to master'
$ws.Range("E3").Value = 'I created maze concept art for level 3'
$ws.Range("E4").Value = 'No significant progress'
$ws.Range("E8").Value = 'Hopefully have the right design down for the help and about pages'
$ws.Range("E12").Value = 'Travel and Time Difference'
$ws.Range("E16").Value = 'It''s easier to get work done when you can communcate with your team '
$ws.Range("E7").Value = 'I will work on my assigned issues'
$ws.Range("E10").Value = 'Not currently'
$ws.Range("E19").Value = 'No'
$ws.Range("E20").Value = 'Not yet'

# --- 3) Row heights grew to fit the newly-wrapped text.
$ws.Rows.Item(2).RowHeight = 166
$ws.Rows.Item(18).RowHeight = 61.5

# --- 4) Leave the selection where the author ended up.
$ws.Range("E20").Select()
